$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 1.38
$ws.Range("B3").Value = 1.46
$ws.Range("G3").Value = 0.53
$ws.Range("C4").Value = 1.44
$ws.Range("F4").Value = 1.08
$ws.Range("D5").Value = 1.3
$ws.Range("F5").Value = 1.02
$ws.Range("D6").Value = 1.56
$ws.Range("E6").Value = 1.35
$ws.Range("G6").Value = 1.04
$ws.Range("C7").Value = 2.35
$ws.Range("F7").Value = 1.46
